# Auto-generated edit script: apply scheduled-runner market price refresh
# to the Leve profit-tracking sheets (currentAveragePrice / LevePrice / LeveProfit
# columns H-N), per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1805.6  # H32: was 1825.6
$ws.Cells.Item(32, 9).Value = 1832.6666  # I32: was 2299
$ws.Cells.Item(32, 10).Value = 1765  # J32: was 1510
$ws.Cells.Item(32, 11).Value = 1832.6666  # K32: was 2299
$ws.Cells.Item(32, 12).Value = 1765  # L32: was 1510
$ws.Cells.Item(32, 13).Value = -1506.6666  # M32: was -1973
$ws.Cells.Item(32, 14).Value = -2417  # N32: was -2162
$ws.Cells.Item(101, 8).Value = 1814.6923  # H101: was 1063.3636
$ws.Cells.Item(101, 9).Value = 2076.7778  # I101: was 971
$ws.Cells.Item(101, 11).Value = 6230.3334  # K101: was 2913
$ws.Cells.Item(101, 13).Value = -4608.3334  # M101: was -1291
$ws.Cells.Item(132, 8).Value = 2412.2856  # H132: was 2383.279
$ws.Cells.Item(132, 9).Value = 2285.125  # I132: was 2297
$ws.Cells.Item(132, 10).Value = 4955.5  # J132: was 3533.6667
$ws.Cells.Item(132, 11).Value = 6855.375  # K132: was 6891
$ws.Cells.Item(132, 12).Value = 14866.5  # L132: was 10601.0001
$ws.Cells.Item(132, 13).Value = -4325.375  # M132: was -4361
$ws.Cells.Item(132, 14).Value = -19926.5  # N132: was -15661.0001
$ws.Cells.Item(138, 8).Value = 3086.3096  # H138: was 3021.5747
$ws.Cells.Item(138, 9).Value = 1989.2  # I138: was 1686.5834
$ws.Cells.Item(138, 10).Value = 3234.5676  # J138: was 3235.1733
$ws.Cells.Item(138, 11).Value = 5967.6  # K138: was 5059.7502
$ws.Cells.Item(138, 12).Value = 9703.702799999999  # L138: was 9705.519899999999
$ws.Cells.Item(138, 13).Value = -827.6000000000004  # M138: was 80.2497999999996
$ws.Cells.Item(138, 14).Value = -19983.7028  # N138: was -19985.5199

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1165.2778  # H2: was 1189.8823
$ws.Cells.Item(2, 9).Value = 1187.7142  # I2: was 1221.6154
$ws.Cells.Item(2, 11).Value = 1187.7142  # K2: was 1221.6154
$ws.Cells.Item(2, 13).Value = -1074.7142  # M2: was -1108.6154
$ws.Cells.Item(32, 8).Value = 12198382  # H32: was 12503354
$ws.Cells.Item(32, 9).Value = 14707974  # I32: was 15153686
$ws.Cells.Item(32, 11).Value = 14707974  # K32: was 15153686
$ws.Cells.Item(32, 13).Value = -14707687  # M32: was -15153399
$ws.Cells.Item(45, 8).Value = 1370.875  # H45: was 1196.8
$ws.Cells.Item(45, 9).Value = 1443  # I45: was 1254.4
$ws.Cells.Item(45, 10).Value = 1298.75  # J45: was 1139.2
$ws.Cells.Item(45, 11).Value = 1443  # K45: was 1254.4
$ws.Cells.Item(45, 12).Value = 1298.75  # L45: was 1139.2
$ws.Cells.Item(45, 13).Value = -1066  # M45: was -877.4000000000001
$ws.Cells.Item(45, 14).Value = -2052.75  # N45: was -1893.2
$ws.Cells.Item(75, 8).Value = 69737  # H75: was 0
$ws.Cells.Item(75, 10).Value = 69737  # J75: was 0
$ws.Cells.Item(75, 12).Value = 69737  # L75: was 0
$ws.Cells.Item(75, 14).Value = -71485  # N75: was None
$ws.Cells.Item(78, 8).Value = 69737  # H78: was 0
$ws.Cells.Item(78, 10).Value = 69737  # J78: was 0
$ws.Cells.Item(78, 12).Value = 209211  # L78: was 0
$ws.Cells.Item(78, 14).Value = -217947  # N78: was None
$ws.Cells.Item(86, 8).Value = 67458  # H86: was 0
$ws.Cells.Item(86, 10).Value = 67458  # J86: was 0
$ws.Cells.Item(86, 12).Value = 67458  # L86: was 0
$ws.Cells.Item(86, 14).Value = -69830  # N86: was None
$ws.Cells.Item(89, 8).Value = 67458  # H89: was 0
$ws.Cells.Item(89, 10).Value = 67458  # J89: was 0
$ws.Cells.Item(89, 12).Value = 202374  # L89: was 0
$ws.Cells.Item(89, 14).Value = -214230  # N89: was None
$ws.Cells.Item(102, 8).Value = 5228.85  # H102: was 5449.8423
$ws.Cells.Item(102, 9).Value = 5267.2104  # I102: was 5502.6113
$ws.Cells.Item(102, 11).Value = 5267.2104  # K102: was 5502.6113
$ws.Cells.Item(102, 13).Value = -3645.2104  # M102: was -3880.6113
$ws.Cells.Item(110, 8).Value = 1910.8  # H110: was 2034.2222
$ws.Cells.Item(110, 9).Value = 1901  # I110: was 2038.625
$ws.Cells.Item(110, 11).Value = 1901  # K110: was 2038.625
$ws.Cells.Item(110, 13).Value = 144  # M110: was 6.375
$ws.Cells.Item(116, 8).Value = 1165.2778  # H116: was 1189.8823
$ws.Cells.Item(116, 9).Value = 1187.7142  # I116: was 1221.6154
$ws.Cells.Item(116, 11).Value = 1187.7142  # K116: was 1221.6154
$ws.Cells.Item(116, 13).Value = 1106.2858  # M116: was 1072.3846
$ws.Cells.Item(122, 8).Value = 2393.0789  # H122: was 2170.558
$ws.Cells.Item(122, 9).Value = 2129.2415  # I122: was 1928.6364
$ws.Cells.Item(122, 10).Value = 3243.2222  # J122: was 2968.9
$ws.Cells.Item(122, 11).Value = 6387.7245  # K122: was 5785.9092
$ws.Cells.Item(122, 12).Value = 9729.6666  # L122: was 8906.700000000001
$ws.Cells.Item(122, 13).Value = -3937.7245  # M122: was -3335.9092
$ws.Cells.Item(122, 14).Value = -14629.6666  # N122: was -13806.7
$ws.Cells.Item(132, 8).Value = 6395.6553  # H132: was 6209.3335
$ws.Cells.Item(132, 9).Value = 1506.7727  # I132: was 1506.6364
$ws.Cells.Item(132, 10).Value = 21760.715  # J132: was 19141.75
$ws.Cells.Item(132, 11).Value = 4520.3181  # K132: was 4519.9092
$ws.Cells.Item(132, 12).Value = 65282.145  # L132: was 57425.25
$ws.Cells.Item(132, 13).Value = -1990.3181  # M132: was -1989.9092
$ws.Cells.Item(132, 14).Value = -70342.145  # N132: was -62485.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1165.2778  # H3: was 1189.8823
$ws.Cells.Item(3, 9).Value = 1187.7142  # I3: was 1221.6154
$ws.Cells.Item(3, 11).Value = 1187.7142  # K3: was 1221.6154
$ws.Cells.Item(3, 13).Value = -1073.7142  # M3: was -1107.6154
$ws.Cells.Item(76, 8).Value = 47513.2  # H76: was 52513.2
$ws.Cells.Item(76, 10).Value = 47513.2  # J76: was 52513.2
$ws.Cells.Item(76, 12).Value = 47513.2  # L76: was 52513.2
$ws.Cells.Item(76, 14).Value = -48143.2  # N76: was -53143.2
$ws.Cells.Item(79, 8).Value = 47513.2  # H79: was 52513.2
$ws.Cells.Item(79, 10).Value = 47513.2  # J79: was 52513.2
$ws.Cells.Item(79, 12).Value = 47513.2  # L79: was 52513.2
$ws.Cells.Item(79, 14).Value = -49697.2  # N79: was -54697.2
$ws.Cells.Item(107, 8).Value = 2513.2778  # H107: was 2602.353
$ws.Cells.Item(107, 9).Value = 2214.2856  # I107: was 2307.7693
$ws.Cells.Item(107, 11).Value = 2214.2856  # K107: was 2307.7693
$ws.Cells.Item(107, 13).Value = -294.2856000000002  # M107: was -387.7692999999999
$ws.Cells.Item(134, 8).Value = 401244.9  # H134: was 358269.1
$ws.Cells.Item(134, 9).Value = 758.63635  # I134: was 684.08
$ws.Cells.Item(134, 11).Value = 2275.90905  # K134: was 2052.24
$ws.Cells.Item(134, 13).Value = 259.0909499999998  # M134: was 482.7599999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1642.1666  # H58: was 1614.75
$ws.Cells.Item(58, 9).Value = 1642.1666  # I58: was 1607.4286
$ws.Cells.Item(58, 10).Value = 0  # J58: was 1666
$ws.Cells.Item(58, 11).Value = 1642.1666  # K58: was 1607.4286
$ws.Cells.Item(58, 12).Value = 0  # L58: was 1666
$ws.Cells.Item(58, 13).Value = -1439.1666  # M58: was -1404.4286
$ws.Cells.Item(58, 14).Value = $null  # N58: clear (was -2072)
$ws.Cells.Item(132, 8).Value = 2165.9285  # H132: was 1954.3334
$ws.Cells.Item(132, 9).Value = 2121.2195  # I132: was 1971.7446
$ws.Cells.Item(132, 10).Value = 3999  # J132: was 1749.75
$ws.Cells.Item(132, 11).Value = 6363.6585  # K132: was 5915.2338
$ws.Cells.Item(132, 12).Value = 11997  # L132: was 5249.25
$ws.Cells.Item(132, 13).Value = -3833.6585  # M132: was -3385.2338
$ws.Cells.Item(132, 14).Value = -17057  # N132: was -10309.25
$ws.Cells.Item(134, 8).Value = 1001738.6  # H134: was 770840.9399999999
$ws.Cells.Item(134, 9).Value = 1001738.6  # I134: was 770840.9399999999
$ws.Cells.Item(134, 11).Value = 3005215.8  # K134: was 2312522.82
$ws.Cells.Item(134, 13).Value = -3002680.8  # M134: was -2309987.82
$ws.Cells.Item(136, 8).Value = 1642.1666  # H136: was 1614.75
$ws.Cells.Item(136, 9).Value = 1642.1666  # I136: was 1607.4286
$ws.Cells.Item(136, 10).Value = 0  # J136: was 1666
$ws.Cells.Item(136, 11).Value = 4926.4998  # K136: was 4822.2858
$ws.Cells.Item(136, 12).Value = 0  # L136: was 4998
$ws.Cells.Item(136, 13).Value = -2376.4998  # M136: was -2272.2858
$ws.Cells.Item(136, 14).Value = $null  # N136: clear (was -10098)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(109, 8).Value = 2934.7778  # H109: was 3127.875
$ws.Cells.Item(109, 9).Value = 3104.75  # I109: was 3616.3333
$ws.Cells.Item(109, 10).Value = 1575  # J109: was 1662.5
$ws.Cells.Item(109, 11).Value = 9314.25  # K109: was 10848.9999
$ws.Cells.Item(109, 12).Value = 4725  # L109: was 4987.5
$ws.Cells.Item(109, 13).Value = -8274.25  # M109: was -9808.999899999999
$ws.Cells.Item(109, 14).Value = -6805  # N109: was -7067.5
$ws.Cells.Item(131, 8).Value = 6196.079  # H131: was 6209.6055
$ws.Cells.Item(131, 10).Value = 5679.7812  # J131: was 5695.8438
$ws.Cells.Item(131, 12).Value = 17039.3436  # L131: was 17087.5314
$ws.Cells.Item(131, 14).Value = -27119.3436  # N131: was -27167.5314
$ws.Cells.Item(132, 8).Value = 2189.8572  # H132: was 2183.842
$ws.Cells.Item(132, 9).Value = 1997.25  # I132: was 1997.5454
$ws.Cells.Item(132, 10).Value = 2446.6667  # J132: was 2440
$ws.Cells.Item(132, 11).Value = 17975.25  # K132: was 17977.9086
$ws.Cells.Item(132, 12).Value = 22020.0003  # L132: was 21960
$ws.Cells.Item(132, 13).Value = -15445.25  # M132: was -15447.9086
$ws.Cells.Item(132, 14).Value = -27080.0003  # N132: was -27020
$ws.Cells.Item(134, 8).Value = 3288.2856  # H134: was 2813.7646
$ws.Cells.Item(134, 9).Value = 924.6  # I134: was 849.53845
$ws.Cells.Item(134, 11).Value = 2773.8  # K134: was 2548.61535
$ws.Cells.Item(134, 13).Value = 2296.2  # M134: was 2521.38465

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(64, 8).Value = 59879  # H64: was 59874.5
$ws.Cells.Item(64, 10).Value = 59879  # J64: was 59874.5
$ws.Cells.Item(64, 12).Value = 59879  # L64: was 59874.5
$ws.Cells.Item(64, 14).Value = -60375  # N64: was -60370.5
$ws.Cells.Item(67, 8).Value = 59879  # H67: was 59874.5
$ws.Cells.Item(67, 10).Value = 59879  # J67: was 59874.5
$ws.Cells.Item(67, 12).Value = 59879  # L67: was 59874.5
$ws.Cells.Item(67, 14).Value = -61595  # N67: was -61590.5
$ws.Cells.Item(132, 8).Value = 125016570  # H132: was 166687790
$ws.Cells.Item(132, 9).Value = 250002180  # I132: was 500001440
$ws.Cells.Item(132, 11).Value = 750006540  # K132: was 1500004320
$ws.Cells.Item(132, 13).Value = -750004010  # M132: was -1500001790
$ws.Cells.Item(133, 8).Value = 0  # H133: was 150000
$ws.Cells.Item(133, 10).Value = 0  # J133: was 150000
$ws.Cells.Item(133, 12).Value = 0  # L133: was 150000
$ws.Cells.Item(133, 14).Value = $null  # N133: clear (was -160120)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 129451.375  # H7: was 147658.72
$ws.Cells.Item(16, 8).Value = 1265.963  # H16: was 1457.8077
$ws.Cells.Item(16, 9).Value = 1119.9131  # I16: was 1213.3334
$ws.Cells.Item(16, 10).Value = 2105.75  # J16: was 2484.6
$ws.Cells.Item(16, 11).Value = 1119.9131  # K16: was 1213.3334
$ws.Cells.Item(16, 12).Value = 2105.75  # L16: was 2484.6
$ws.Cells.Item(16, 13).Value = -949.9131  # M16: was -1043.3334
$ws.Cells.Item(16, 14).Value = -2445.75  # N16: was -2824.6
$ws.Cells.Item(123, 8).Value = 0  # H123: was 46000
$ws.Cells.Item(123, 10).Value = 0  # J123: was 46000
$ws.Cells.Item(123, 12).Value = 0  # L123: was 46000
$ws.Cells.Item(123, 14).Value = $null  # N123: clear (was -55800)
$ws.Cells.Item(125, 8).Value = 103750  # H125: was 92698
$ws.Cells.Item(125, 10).Value = 103750  # J125: was 92698
$ws.Cells.Item(125, 12).Value = 103750  # L125: was 92698
$ws.Cells.Item(125, 14).Value = -113590  # N125: was -102538
$ws.Cells.Item(126, 8).Value = 129451.375  # H126: was 147658.72
$ws.Cells.Item(132, 8).Value = 836520.3  # H132: was 2503001.5
$ws.Cells.Item(132, 9).Value = 3323.8  # I132: was 3500
$ws.Cells.Item(132, 11).Value = 9971.400000000001  # K132: was 10500
$ws.Cells.Item(132, 13).Value = -7441.400000000001  # M132: was -7970

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 65000  # H46: was 63333.332
$ws.Cells.Item(46, 10).Value = 65000  # J46: was 63333.332
$ws.Cells.Item(46, 12).Value = 65000  # L46: was 63333.332
$ws.Cells.Item(46, 14).Value = -65462  # N46: was -63795.332
$ws.Cells.Item(88, 8).Value = 0  # H88: was 40000
$ws.Cells.Item(88, 10).Value = 0  # J88: was 40000
$ws.Cells.Item(88, 12).Value = 0  # L88: was 40000
$ws.Cells.Item(88, 14).Value = $null  # N88: clear (was -40812)
$ws.Cells.Item(91, 8).Value = 0  # H91: was 40000
$ws.Cells.Item(91, 10).Value = 0  # J91: was 40000
$ws.Cells.Item(91, 12).Value = 0  # L91: was 40000
$ws.Cells.Item(91, 14).Value = $null  # N91: clear (was -42808)
$ws.Cells.Item(132, 8).Value = 2312.0356  # H132: was 2339.8215
$ws.Cells.Item(132, 9).Value = 2375.087  # I132: was 2472.9546
$ws.Cells.Item(132, 10).Value = 2022  # J132: was 1851.6666
$ws.Cells.Item(132, 11).Value = 7125.261  # K132: was 7418.8638
$ws.Cells.Item(132, 12).Value = 6066  # L132: was 5554.9998
$ws.Cells.Item(132, 13).Value = -4595.261  # M132: was -4888.8638
$ws.Cells.Item(132, 14).Value = -11126  # N132: was -10614.9998
$ws.Cells.Item(134, 8).Value = 65000  # H134: was 63333.332
$ws.Cells.Item(134, 10).Value = 65000  # J134: was 63333.332
$ws.Cells.Item(134, 12).Value = 195000  # L134: was 189999.996
$ws.Cells.Item(134, 14).Value = -200070  # N134: was -195069.996

